$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "trainingaudio/18_popata2.wav"
$ws.Range("B2").Value = "pngimages/18_donut.png"

$ws.Range("A3").Value = "trainingaudio/21_papika1.wav"
$ws.Range("B3").Value = "pngimages/21_cheese.png"

$ws.Range("A4").Value = "trainingaudio/23_patoko1.wav"
$ws.Range("B4").Value = "pngimages/23_lemon.png"

$ws.Range("A5").Value = "trainingaudio/06_titoka3.wav"
$ws.Range("B5").Value = "pngimages/06_tent.png"

$ws.Range("A6").Value = "trainingaudio/04_kitoti2.wav"
$ws.Range("B6").Value = "pngimages/04_ladder.png"

$ws.Range("A7").Value = "trainingaudio/05_titopo2.wav"
$ws.Range("B7").Value = "pngimages/05_megaphone.png"
